$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 5 (Natalie's - Strawberry Lemonade) ---
# Quantity changes from 2 -> 1
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "1"

# Total Cost changes from 18.50 -> 9.25 (Quantity 1 * Cost Per 9.25)
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "9.25"

# --- Add new row 10 for Natalie's - Honey Tangerine ---
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "004061"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "Natalie's - Honey Tangerine"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "1"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "14.00"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "14.00"
